$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value2 = 4666.6665
$ws.Range("I21").Value2 = 6000
$ws.Range("K21").Value2 = 6000
$ws.Range("M21").Value2 = -5532
$ws.Range("H23").Value2 = 4666.6665
$ws.Range("I23").Value2 = 6000
$ws.Range("K23").Value2 = 6000
$ws.Range("M23").Value2 = -5766
$ws.Range("H29").Value2 = 1855.5555
$ws.Range("J29").Value2 = 2700
$ws.Range("L29").Value2 = 8100
$ws.Range("N29").Value2 = -8662
$ws.Range("H38").Value2 = 1608
$ws.Range("J38").Value2 = 4225
$ws.Range("L38").Value2 = 12675
$ws.Range("N38").Value2 = -13419
$ws.Range("H58").Value2 = 1060.5
$ws.Range("I58").Value2 = 483
$ws.Range("J58").Value2 = 2504.25
$ws.Range("K58").Value2 = 1449
$ws.Range("L58").Value2 = 7512.75
$ws.Range("M58").Value2 = -1299
$ws.Range("N58").Value2 = -7812.75
$ws.Range("H100").Value2 = 901.6
$ws.Range("I100").Value2 = 632.61536
$ws.Range("J100").Value2 = 2650
$ws.Range("K100").Value2 = 632.61536
$ws.Range("L100").Value2 = 2650
$ws.Range("M100").Value2 = -91.61536000000001
$ws.Range("N100").Value2 = -3732
$ws.Range("H132").Value2 = 9811493
$ws.Range("I132").Value2 = 13895918
$ws.Range("J132").Value2 = 8871.700000000001
$ws.Range("K132").Value2 = 41687754
$ws.Range("L132").Value2 = 26615.1
$ws.Range("M132").Value2 = -41685224
$ws.Range("N132").Value2 = -31675.1
$ws.Range("H135").Value2 = 43479644
$ws.Range("J135").Value2 = 166670700
$ws.Range("L135").Value2 = 1500036300
$ws.Range("N135").Value2 = -1500041370
$ws.Range("H137").Value2 = 1415.5834
$ws.Range("I137").Value2 = 1026.6842
$ws.Range("K137").Value2 = 3080.0526
$ws.Range("M137").Value2 = -530.0526
$ws.Range("H141").Value2 = 2403.5
$ws.Range("I141").Value2 = 2403.5
$ws.Range("K141").Value2 = 7210.5
$ws.Range("M141").Value2 = -2030.5

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value2 = 4910.2925
$ws.Range("I32").Value2 = 4809.1934
$ws.Range("K32").Value2 = 4809.1934
$ws.Range("M32").Value2 = -4522.1934
$ws.Range("H45").Value2 = 1651.4706
$ws.Range("I45").Value2 = 1606.3334
$ws.Range("K45").Value2 = 1606.3334
$ws.Range("M45").Value2 = -1229.3334
$ws.Range("H54").Value2 = 0
$ws.Range("J54").Value2 = 0
$ws.Range("L54").Value2 = 0
$ws.Range("N54").ClearContents()
$ws.Range("H61").Value2 = 71429896
$ws.Range("I61").Value2 = 90910140
$ws.Range("K61").Value2 = 90910140
$ws.Range("M61").Value2 = -90909928
$ws.Range("H63").Value2 = 22224056
$ws.Range("J63").Value2 = 71430780
$ws.Range("L63").Value2 = 71430780
$ws.Range("N63").Value2 = -71432152
$ws.Range("H66").Value2 = 22224056
$ws.Range("J66").Value2 = 71430780
$ws.Range("L66").Value2 = 357153900
$ws.Range("N66").Value2 = -357160764
$ws.Range("H110").Value2 = 1182.8889
$ws.Range("I110").Value2 = 773.3077
$ws.Range("J110").Value2 = 2247.8
$ws.Range("K110").Value2 = 773.3077
$ws.Range("L110").Value2 = 2247.8
$ws.Range("M110").Value2 = 1271.6923
$ws.Range("N110").Value2 = -6337.8
$ws.Range("H122").Value2 = 1788.75
$ws.Range("I122").Value2 = 1345.1666
$ws.Range("J122").Value2 = 3119.5
$ws.Range("K122").Value2 = 4035.4998
$ws.Range("L122").Value2 = 9358.5
$ws.Range("M122").Value2 = -1585.4998
$ws.Range("N122").Value2 = -14258.5
$ws.Range("H136").Value2 = 71429896
$ws.Range("I136").Value2 = 90910140
$ws.Range("K136").Value2 = 272730420
$ws.Range("M136").Value2 = -272727870

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value2 = 2831.45
$ws.Range("I86").Value2 = 2853.3076
$ws.Range("J86").Value2 = 2790.8572
$ws.Range("K86").Value2 = 2853.3076
$ws.Range("L86").Value2 = 2790.8572
$ws.Range("M86").Value2 = -1730.3076
$ws.Range("N86").Value2 = -5036.8572
$ws.Range("H89").Value2 = 2831.45
$ws.Range("I89").Value2 = 2853.3076
$ws.Range("J89").Value2 = 2790.8572
$ws.Range("K89").Value2 = 14266.538
$ws.Range("L89").Value2 = 13954.286
$ws.Range("M89").Value2 = -8650.538
$ws.Range("N89").Value2 = -25186.286
$ws.Range("H94").Value2 = 19231124
$ws.Range("I94").Value2 = 31250302
$ws.Range("J94").Value2 = 440
$ws.Range("K94").Value2 = 31250302
$ws.Range("L94").Value2 = 440
$ws.Range("M94").Value2 = -31249851
$ws.Range("N94").Value2 = -1342
$ws.Range("H99").Value2 = 62501172
$ws.Range("I99").Value2 = 66667804
$ws.Range("J99").Value2 = 1700
$ws.Range("K99").Value2 = 66667804
$ws.Range("L99").Value2 = 1700
$ws.Range("M99").Value2 = -66666306
$ws.Range("N99").Value2 = -4696
$ws.Range("H107").Value2 = 1532.3334
$ws.Range("I107").Value2 = 779.75
$ws.Range("J107").Value2 = 2134.4
$ws.Range("K107").Value2 = 779.75
$ws.Range("L107").Value2 = 2134.4
$ws.Range("M107").Value2 = 1140.25
$ws.Range("N107").Value2 = -5974.4
$ws.Range("H134").Value2 = 1388.2
$ws.Range("I134").Value2 = 1075.091
$ws.Range("K134").Value2 = 3225.273
$ws.Range("M134").Value2 = -690.2729999999997

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value2 = 457.25
$ws.Range("I10").Value2 = 457.25
$ws.Range("K10").Value2 = 457.25
$ws.Range("M10").Value2 = -318.25
$ws.Range("H19").Value2 = 260
$ws.Range("J19").Value2 = 403.33334
$ws.Range("L19").Value2 = 403.33334
$ws.Range("N19").Value2 = -743.33334
$ws.Range("H24").Value2 = 260
$ws.Range("J24").Value2 = 403.33334
$ws.Range("L24").Value2 = 403.33334
$ws.Range("N24").Value2 = -743.33334
$ws.Range("H31").Value2 = 1300.3043
$ws.Range("I31").Value2 = 1227.6136
$ws.Range("J31").Value2 = 2899.5
$ws.Range("K31").Value2 = 1227.6136
$ws.Range("L31").Value2 = 2899.5
$ws.Range("M31").Value2 = -932.6135999999999
$ws.Range("N31").Value2 = -3489.5
$ws.Range("H34").Value2 = 1300.3043
$ws.Range("I34").Value2 = 1227.6136
$ws.Range("J34").Value2 = 2899.5
$ws.Range("K34").Value2 = 1227.6136
$ws.Range("L34").Value2 = 2899.5
$ws.Range("M34").Value2 = -1025.6136
$ws.Range("N34").Value2 = -3303.5
$ws.Range("H132").Value2 = 1850.8462
$ws.Range("I132").Value2 = 1574.1666
$ws.Range("K132").Value2 = 4722.4998
$ws.Range("M132").Value2 = -2192.4998
$ws.Range("H134").Value2 = 14287264
$ws.Range("I134").Value2 = 1590.4546
$ws.Range("K134").Value2 = 4771.3638
$ws.Range("M134").Value2 = -2236.3638

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value2 = 631.6667
$ws.Range("I32").Value2 = 631.6667
$ws.Range("J32").Value2 = 0
$ws.Range("K32").Value2 = 1895.0001
$ws.Range("L32").Value2 = 0
$ws.Range("M32").Value2 = -1612.0001
$ws.Range("N32").ClearContents()
$ws.Range("H60").Value2 = 1717.5834
$ws.Range("I60").Value2 = 601.6667
$ws.Range("J60").Value2 = 2089.5557
$ws.Range("K60").Value2 = 1805.0001
$ws.Range("L60").Value2 = 6268.6671
$ws.Range("M60").Value2 = -1554.0001
$ws.Range("N60").Value2 = -6770.6671
$ws.Range("H87").Value2 = 4400
$ws.Range("I87").Value2 = 800
$ws.Range("J87").Value2 = 8000
$ws.Range("K87").Value2 = 2400
$ws.Range("L87").Value2 = 24000
$ws.Range("M87").Value2 = -1152
$ws.Range("N87").Value2 = -26496
$ws.Range("H90").Value2 = 4400
$ws.Range("I90").Value2 = 800
$ws.Range("J90").Value2 = 8000
$ws.Range("K90").Value2 = 7200
$ws.Range("L90").Value2 = 72000
$ws.Range("M90").Value2 = -960
$ws.Range("N90").Value2 = -84480

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value2 = 2926.5454
$ws.Range("I122").Value2 = 3268.1428
$ws.Range("J122").Value2 = 2328.75
$ws.Range("K122").Value2 = 9804.428400000001
$ws.Range("L122").Value2 = 6986.25
$ws.Range("M122").Value2 = -7354.428400000001
$ws.Range("N122").Value2 = -11886.25
$ws.Range("H126").Value2 = 2246.6667
$ws.Range("I126").Value2 = 1908
$ws.Range("J126").Value2 = 2670
$ws.Range("K126").Value2 = 5724
$ws.Range("L126").Value2 = 8010
$ws.Range("M126").Value2 = -3254
$ws.Range("N126").Value2 = -12950
$ws.Range("H132").Value2 = 2692.0293
$ws.Range("I132").Value2 = 2485.9048
$ws.Range("K132").Value2 = 7457.714399999999
$ws.Range("M132").Value2 = -4927.714399999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value2 = 3767.6667
$ws.Range("I40").Value2 = 2107.75
$ws.Range("K40").Value2 = 2107.75
$ws.Range("M40").Value2 = -1971.75
$ws.Range("H61").Value2 = 940
$ws.Range("I61").Value2 = 905.7143
$ws.Range("J61").Value2 = 1100
$ws.Range("K61").Value2 = 905.7143
$ws.Range("L61").Value2 = 1100
$ws.Range("M61").Value2 = -703.7143
$ws.Range("N61").Value2 = -1504
$ws.Range("H100").Value2 = 1140.1666
$ws.Range("I100").Value2 = 985.25
$ws.Range("K100").Value2 = 985.25
$ws.Range("M100").Value2 = -444.25
$ws.Range("H113").Value2 = 940
$ws.Range("I113").Value2 = 905.7143
$ws.Range("J113").Value2 = 1100
$ws.Range("K113").Value2 = 905.7143
$ws.Range("L113").Value2 = 1100
$ws.Range("M113").Value2 = 1264.2857
$ws.Range("N113").Value2 = -5440
$ws.Range("H133").Value2 = 47306.5
$ws.Range("J133").Value2 = 47306.5
$ws.Range("L133").Value2 = 47306.5
$ws.Range("N133").Value2 = -52366.5
$ws.Range("H139").Value2 = 70715
$ws.Range("J139").Value2 = 70715
$ws.Range("L139").Value2 = 70715
$ws.Range("N139").Value2 = -80995

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value2 = 304.5
$ws.Range("I107").Value2 = 258.33334
$ws.Range("J107").Value2 = 387.6
$ws.Range("K107").Value2 = 775.0000200000001
$ws.Range("L107").Value2 = 1162.8
$ws.Range("M107").Value2 = 1144.99998
$ws.Range("N107").Value2 = -5002.8
$ws.Range("H126").Value2 = 76928390
$ws.Range("I126").Value2 = 111114456
$ws.Range("K126").Value2 = 333343368
$ws.Range("M126").Value2 = -333340898
$ws.Range("H132").Value2 = 2203.348
$ws.Range("I132").Value2 = 2084
$ws.Range("K132").Value2 = 6252
$ws.Range("M132").Value2 = -3722
$ws.Range("H136").Value2 = 1076.7391
$ws.Range("J136").Value2 = 1516.2727
$ws.Range("L136").Value2 = 4548.8181
$ws.Range("N136").Value2 = -9648.8181
